$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (old rows 4 and 5); the remaining rows 2-3
# will be updated in place with recalculated TPM-based values below.
$ws.Rows("4:5").Delete()

# Row 2: Sending cluster = ECs, Target cluster = MuSCs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Agrp"
$ws.Range("C2").Value = "Mc5r"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3518616666666667
$ws.Range("H2").Value = 1.055585
$ws.Range("I2").Value = 0.5958054833396739
$ws.Range("J2").Value = 0.5958054833396738
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02331333333333334
$ws.Range("N2").Value = 0.06994
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.008203068322222224
$ws.Range("R2").Value = 0.0738276149
$ws.Range("S2").Value = 0.5958054833396739
$ws.Range("T2").Value = 0.5958054833396738

# Row 3: Sending cluster = MuSCs, Target cluster = MuSCs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Agrp"
$ws.Range("C3").Value = "Mc5r"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.238703
$ws.Range("H3").Value = 0.716109
$ws.Range("I3").Value = 0.4041945166603262
$ws.Range("J3").Value = 0.4041945166603262
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.02331333333333334
$ws.Range("N3").Value = 0.06994
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.005564962606666667
$ws.Range("R3").Value = 0.05008466346
$ws.Range("S3").Value = 0.4041945166603262
$ws.Range("T3").Value = 0.4041945166603262
